# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch count)
# for column G, rows 2-34 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 5
    4  = 4
    5  = 9
    6  = 7
    7  = 5
    8  = 7
    9  = 4
    10 = 11
    11 = 12
    12 = 7
    13 = 8
    14 = 5
    15 = 7
    16 = 6
    17 = 7
    18 = 10
    19 = 5
    20 = 11
    21 = 9
    22 = 12
    23 = 8
    24 = 5
    25 = 4
    26 = 5
    27 = 9
    28 = 5
    29 = 3
    30 = 0
    31 = 4
    32 = 5
    33 = 5
    34 = 5
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
